# Applies the "Updated cryptos list" data refresh to sheet1.
# Updates price (column D) and 1h volume-change (column E) values for many rows,
# and swaps the Cronos/FraxShare row content (rows 45/46) per the source diff.
#
# Price cells are forced to Text format before assignment (and the format is
# reset back to the default "Normal" style afterward) so that values such as
# "229.21" or "0.0610" are stored as literal text -- matching the original
# workbook -- instead of being auto-converted to floating point numbers by Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "38.704.24"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.04%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.092.71"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.24%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "229.21"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.27%  "
$ws.Range("E6").Value = "  +0.48%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "61.27"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.63%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("E9").Value = "  +1.54%  "
$ws.Range("E10").Value = "  +0.73%  "
$ws.Range("E11").Value = "  -0.51%  "
$ws.Range("E12").Value = "  +4.94%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.404.18"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "22.13"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.68%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.818"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +6.01%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.48"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.05%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.105.49"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.56%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "38.697.06"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.14%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "71.68"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.49%  "
$ws.Range("E20").Value = "  +1.07%  "
$ws.Range("E21").Value = "  +0.55%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "227.28"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.59%  "
$ws.Range("E24").Value = "  -1.52%  "
$ws.Range("E25").Value = "  +0.70%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.58"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.89%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "171.30"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.01%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.139"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +6.89%  "
$ws.Range("E29").Value = "  +9.95%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.26"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.48%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.47"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.58%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.120"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.44%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.52"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.30%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.69"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.02%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0610"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.53%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.53"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.95%  "
$ws.Range("E37").Value = "  -0.10%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.60"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.84%  "
$ws.Range("E39").Value = "  -0.05%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.06"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.21%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0230"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.27%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "101.11"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.10%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.533.98"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.55%  "
$ws.Range("E44").Value = "  -0.57%  "
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.79"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +8.33%  "
$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0916"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.18%  "
$ws.Range("E47").Value = "  +1.05%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.10"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.10%  "
$ws.Range("E49").Value = "  +1.09%  "
$ws.Range("E50").Value = "  -0.70%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.290.24"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.34%  "
